$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 9;  Date = 42613.758136574077; B = 81 },
    @{ Row = 10; Date = 42613.885868055557; B = 7 },
    @{ Row = 11; Date = 42614.884247685186; B = 50 },
    @{ Row = 12; Date = 42615.884733796294; B = 10 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Date
    $ws.Cells.Item($row, 2).Value = $r.B
    for ($col = 3; $col -le 13; $col++) {
        $ws.Cells.Item($row, $col).Value = 0
    }
    $ws.Cells.Item($row, 14).Value = "Random"
}
